$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176252841949463
$ws.Range("B1").Value = 2.411530256271362
$ws.Range("D1").Value = 2.339471340179443
$ws.Range("E1").Value = 1.201716423034668
